# Updating exec time table: insert a new "ignore-alias /clean" column between
# the existing Benchmark (A) and opt/clean (B) columns, pushing the old
# opt/clean data into column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the existing "opt/clean" column (with its data + formatting)
#        from B into C, then mirror C's formatting back onto B so the new
#        column inherits the same borders/number formats. ---
$ws.Range("B1:B32").Copy($ws.Range("C1:C32"))
$ws.Range("C1:C32").Copy($ws.Range("B1:B32"))

# --- 2. Header text ---
$ws.Range("B1").Value2 = "ignore-alias /clean"
$ws.Range("C1").Value2 = "opt /clean"

# --- 3. New "ignore-alias /clean" data values for rows 2-31 ---
$newVals = @(
    0.99684918725841232,
    0.99786251461118336,
    0.99885250667016012,
    1.0115700171216659,
    2.0359277488759853,
    0.93816242191146082,
    0.99078206538850777,
    0.90427142189875376,
    0.97385654844127834,
    0.43822490563299604,
    1.1361512416553816,
    0.28315643691924919,
    0.99421835594677654,
    0.97185107990721742,
    0.94683263407344598,
    0.9856877599300119,
    1.0015622033817697,
    1.0061330375189017,
    0.97086344982154482,
    1.0041683420844889,
    0.99792136700503431,
    1.0025008857060866,
    1.0011620401056651,
    2.1359834390407868,
    1.0012426025070573,
    1.0011265099134947,
    1.0028817636109317,
    1.7536201903267143,
    1.7732759876044697,
    0.99549626378357969
)

$arr = New-Object 'object[,]' 30,1
for ($i = 0; $i -lt 30; $i++) {
    $arr[$i,0] = $newVals[$i]
}
$ws.Range("B2:B31").Value2 = $arr

Write-Host "done"
